$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 990.5413814561294
$ws.Range("F2").Value = -0.9458618543870534

$ws.Range("E4").Value = 890.4766276913342
$ws.Range("F4").Value = -10.95233723086658

$ws.Range("E5").Value = 888.6058483520386
$ws.Range("F5").Value = -11.13941516479614

$ws.Range("E8").Value = 1024.494620930697
$ws.Range("F8").Value = 2.449462093069699

$ws.Range("E10").Value = 991.3163165158546
$ws.Range("F10").Value = -0.8683683484145388
